# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to match the newly scraped data at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F on both affected sheets.
$updates = @{
    2  = 344
    4  = 10575
    7  = 93
    8  = 1300
    9  = 8081
    10 = 25
    11 = 460
    12 = 2
    14 = 133
    15 = 3244
    16 = 39
    18 = 730
    20 = 1051
    22 = 95
    23 = 1693
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
